$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 27; this shifts the existing rows 27-35 down to 28-36,
# preserving their content and formatting (including the date style on column D).
$ws.Rows(27).Insert()

# Populate the newly inserted row 27 with the new record's data.
$ws.Range("A27").Value = 3
$ws.Range("B27").Value = "Femacal de La Calera"
$ws.Range("C27").Value = "Coquimbo"
$ws.Range("D27").Value = 44511
$ws.Range("D27").NumberFormat = $ws.Range("D28").NumberFormat
$ws.Range("E27").Value = 5
$ws.Range("F27").Value = 100112022
$ws.Range("G27").Value = "Arveja Verde"
$ws.Range("H27").Value = "Perfection"
$ws.Range("I27").Value = "Primera"
$ws.Range("J27").Value = 73
$ws.Range("K27").Value = 16000
$ws.Range("L27").Value = 17000
$ws.Range("M27").Value = 16479
$ws.Range("N27").Value = "$/saco 25 kilos"
$ws.Range("O27").Value = "Provincia de Limarí"
$ws.Range("P27").Value = 659
$ws.Range("Q27").Value = 25
$ws.Range("R27").Value = "Hortaliza"
